# Workbook edit: convert column C ("Bénéfice (après 2 ans") from a
# percentage-style fraction (0.05, 0.1, ...) to a plain integer value
# (5, 10, ...) and drop the explicit number-format style that column
# carried, per the commit "pb int vs float dans les conversions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New values for C2:C21 (old value * 100, stored as plain numbers) ---
$newValues = @{
    2  = 5
    3  = 10
    4  = 15
    5  = 20
    6  = 17
    7  = 25
    8  = 7
    9  = 11
    10 = 13
    11 = 27
    12 = 17
    13 = 9
    14 = 23
    15 = 1
    16 = 3
    17 = 8
    18 = 12
    19 = 14
    20 = 21
    21 = 18
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}

# --- Drop the extra/no-longer-used style (numFmtId 0, applyNumberFormat)
# that used to be applied to C2:C21 (s="2"), reverting those cells back
# to the default (unstyled) cellXfs entry. C1 keeps its bold header
# style, so restore it explicitly after resetting the column.
$ws.Columns.Item(3).Style = "Normal"
$ws.Range("C1").Font.Bold = $true

# --- Update the current view state: scrolled so row 14 is at the top,
# with D21 as the active selected cell.
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D21").Select()
